$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.859.76"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.839.24"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "231.42"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "39.77"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "0.0685"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "0.0987"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "2.106.19"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "11.45"
$ws.Range("E13").Value = "  +3.71%  "
$ws.Range("D14").Value = "1.837.30"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "4.64"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "34.876.34"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "69.82"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "240.27"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "4.69"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "171.10"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "17.42"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.124"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "0.0552"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("E36").Value = "  +11.47%  "
$ws.Range("D37").Value = "0.693"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("E38").Value = "  +6.95%  "
$ws.Range("D39").Value = "91.19"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").Value = "1.340.39"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "14.74"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Value = "2.27"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "6.28"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D48").Value = "2.019.29"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "0.0664"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("E51").Value = "  +16.24%  "
